$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray #NUM! error literals that were left in column F for
# these rows (the cells are deleted outright, not just blanked).
$errorCells = @("F30","F31","F46","F47","F48","F49","F56","F57","F58","F59")
foreach ($cellRef in $errorCells) {
    $ws.Range($cellRef).ClearContents()
}

# Update the saved selection on the sheet from H14 to A2.
[void]$ws.Range("A2").Select()
